# Reorder the sheets so that "review_info" becomes the first sheet and
# "hotel_info" becomes the second sheet.
$wb = $excel.ActiveWorkbook
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

# Add a new "State" column to hotel_info, right after "Hotel_Name" (i.e.
# before "City"), and populate the single data row with "Louisiana".
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns("C:C").Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"
